$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 contains the "Dummy Product 14" test product entry; bump it to 15.
$ws.Range("A2").Value = "Test - Dummy Product 15"
$ws.Range("B2").Value = "Dummy Product 15"
$ws.Range("C2").Value = "TEST - Dummy 15"
